# Update the contract row for order_id=5 (Weiler / Sandstein):
#  - amount: 1 -> 10
#  - phase:  "Planung" -> "Sprengung" (new phase/status string)
#  - price:  75 -> 750
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D6").Value = 10
$ws.Range("F6").Value = "Sprengung"
$ws.Range("G6").Value = 750
